$wb = $excel.ActiveWorkbook

# --- Create Sheet2 and Sheet3, positioned after Sheet1 ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Sheet3"

# --- Populate Sheet1 ---
$ws1.Range("A1").Value = "Greeting"
$ws1.Range("B1").Value = "Communication"
$ws1.Range("C1").Value = "ID"

$ws1.Range("A2").Value = "hello"
$ws1.Range("B2").Value = "text"
$ws1.Range("C2").Value = 1

$ws1.Range("A3").Value = "bye"
$ws1.Range("B3").Value = "message"
$ws1.Range("C3").Value = 143

$ws1.Range("A4").Value = "solo"
$ws1.Range("B4").Value = "call"
$ws1.Range("C4").Value = 7689

# --- Populate Sheet2 ---
$ws2.Range("A1").Value = "Greeting"
$ws2.Range("B1").Value = "Communication"
$ws2.Range("C1").Value = "ID"

$ws2.Range("A2").Value = "howdy sheet 2"
$ws2.Range("B2").Value = "sheet 2 row 1"
$ws2.Range("C2").Value = 45

$ws2.Range("A3").Value = "later sheet 2"
$ws2.Range("B3").Value = "sheet 2 row 2"
$ws2.Range("C3").Value = 87

$ws2.Range("A4").Value = "what up sheet 2"
$ws2.Range("B4").Value = "sheet 2 row 3"
$ws2.Range("C4").Value = 565

$ws2.Range("A5").Value = "sup sheet 2"
$ws2.Range("B5").Value = "sheet 2 row 4"
$ws2.Range("C5").Value = 0

# --- Populate Sheet3 ---
$ws3.Range("A1").Value = "Greeting"
$ws3.Range("B1").Value = "Communication"
$ws3.Range("C1").Value = "ID"

$ws3.Range("A2").Value = "hi sheet 3"
$ws3.Range("B2").Value = "sheet 3 text to prove the first row"
$ws3.Range("C2").Value = 75867

$ws3.Range("A3").Value = "bye sheet 3"
$ws3.Range("B3").Value = "sheet 3 text to prove the second row"
$ws3.Range("C3").Value = 970600

# --- Column widths (best fit) on Sheet1 and Sheet3, column B ---
$ws1.Columns.Item(2).AutoFit()
$ws3.Columns.Item(2).AutoFit()

# --- Selections matching the authored view state ---
$null = $ws1.Range("C10").Select()
$null = $ws3.Range("B7").Select()
$null = $ws2.Range("E7").Select()

# --- Make Sheet2 the active/selected sheet (select is applied last) ---
$null = $ws2.Activate()
